# Apply updated crypto price/volume data as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" cells contain plain decimal-looking text (e.g. "219.73") that
# Excel would otherwise auto-convert to a floating point number, losing the
# exact original text/precision. Force those specific cells to Text format
# first so the values are stored verbatim as strings, matching the source data.
$textFormatCells = @(
    'D4', 'D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D13', 'D15', 'D17',
    'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28',
    'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38',
    'D39', 'D40', 'D41', 'D42', 'D44', 'D47', 'D48', 'D49', 'D50', 'D51'
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update each changed cell with its new value.
$ws.Range('D2').Value = '26.278.65'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '1.666.37'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = '219.73'
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('E6').Value = '  -1.63%  '
$ws.Range('D7').Value = '1.006'
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').Value = '0.2666'
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('D9').Value = '0.06338'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').Value = '21.03'
$ws.Range('E10').Value = '  -3.41%  '
$ws.Range('D11').Value = '0.07741'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '1.663.63'
$ws.Range('E12').Value = '  -1.20%  '
$ws.Range('D13').Value = '4.430'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').Value = '1.895.33'
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('D15').Value = '0.5470'
$ws.Range('E15').Value = '  -1.92%  '
$ws.Range('D16').Value = '0.0₅8284'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '64.91'
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('D18').Value = '26.308.30'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = '1.006'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').Value = '4.658'
$ws.Range('E20').Value = '  -2.11%  '
$ws.Range('D21').Value = '194.92'
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('D22').Value = '10.17'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('D23').Value = '6.091'
$ws.Range('E23').Value = '  -4.11%  '
$ws.Range('D24').Value = '1.008'
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('D25').Value = '139.56'
$ws.Range('E25').Value = '  -2.42%  '
$ws.Range('D26').Value = '0.1242'
$ws.Range('E26').Value = '  -3.29%  '
$ws.Range('D27').Value = '7.222'
$ws.Range('E27').Value = '  -3.08%  '
$ws.Range('D28').Value = '16.14'
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').Value = '1.418'
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '0.06108'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('D31').Value = '1.285'
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('D32').Value = '3.594'
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('D33').Value = '3.288'
$ws.Range('E33').Value = '  -5.02%  '
$ws.Range('D34').Value = '1.639'
$ws.Range('E34').Value = '  -2.55%  '
$ws.Range('D35').Value = '0.9804'
$ws.Range('E35').Value = '  -3.12%  '
$ws.Range('D36').Value = '2.426'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').Value = '2.787'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').Value = '0.5779'
$ws.Range('E38').Value = '  -3.97%  '
$ws.Range('D39').Value = '6.066'
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('D40').Value = '0.01601'
$ws.Range('E40').Value = '  -2.61%  '
$ws.Range('D41').Value = '0.8628'
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').Value = '1.005'
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('D43').Value = '1.028.06'
$ws.Range('E43').Value = '  -4.91%  '
$ws.Range('D44').Value = '100.16'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').Value = '1.809.93'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '57.85'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('D48').Value = '1.011'
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('D49').Value = '8.112'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('D50').Value = '0.05194'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').Value = '1.482'
$ws.Range('E51').Value = '  +1.25%  '
